$wb = $excel.ActiveWorkbook

# --- VF-2 sheet: rows 38-43 (timer fields) change datatype from "str" to "float" ---
$wsVF2 = $wb.Worksheets.Item("VF-2")
$wsVF2.Range("C38:C43").Value = "float"

# --- UMC-500 sheet: same datatype change for rows 38-43 ---
$wsUMC = $wb.Worksheets.Item("UMC-500")
$wsUMC.Range("C38:C43").Value = "float"

# --- ST-10 sheet: rows 28-33 (timer fields) change datatype from "str" to "float" ---
$wsST10 = $wb.Worksheets.Item("ST-10")
$wsST10.Range("C28:C33").Value = "float"

# ST-10: drop the "Programmable coolant position" row (row 40), shifting the
# "Coolant level" row (41) up, and give it the "float" datatype it lacked before.
$wsST10.Rows.Item(40).Delete()
$wsST10.Range("C40").Value = "float"

# ST-10: page orientation set to portrait (adds <pageSetup>)
$wsST10.PageSetup.Orientation = 1

# --- Selections / active sheet ---
# Set each sheet's own selection first (selecting a range activates its
# sheet), then activate UMC-500 last so it ends up as the active tab.
$wsVF2.Range("C37:C43").Select()
$wsST10.Range("A40:XFD40").Select()
$wsUMC.Range("C37:C43").Select()
$wsUMC.Activate()
